$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.774.01"
$ws.Range("E2").Value = "  +1.54%  "

$ws.Range("D3").Value = "1.884.94"
$ws.Range("E3").Value = "  +1.53%  "

$ws.Range("E4").Value = "  +0.39%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "333.55"
$ws.Range("E5").Value = "  +1.82%  "

$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4721"
$ws.Range("E7").Value = "  +3.77%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3938"
$ws.Range("E8").Value = "  +0.81%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.46"
$ws.Range("E9").Value = "  -0.75%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08068"
$ws.Range("E10").Value = "  +1.89%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.028"
$ws.Range("E11").Value = "  +1.56%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.27"
$ws.Range("E12").Value = "  +3.93%  "

$ws.Range("D13").Value = "1.893.62"
$ws.Range("E13").Value = "  +1.51%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.989"
$ws.Range("E14").Value = "  +1.09%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.147"
$ws.Range("E15").Value = "  -0.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.009"
$ws.Range("E16").Value = "  +0.68%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.33"
$ws.Range("E17").Value = "  +1.37%  "

$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001049"
$ws.Range("E18").Value = "  +1.98%  "

$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06709"
$ws.Range("E19").Value = "  +0.88%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.38"
$ws.Range("E20").Value = "  +1.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.005"
$ws.Range("E21").Value = "  +0.32%  "

$ws.Range("B22").Value = "WrappedBTC"
$ws.Range("C22").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D22").Value = "27.786.19"
$ws.Range("E22").Value = "  +1.58%  "

$ws.Range("B23").Value = "Uniswap"
$ws.Range("C23").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.547"
$ws.Range("E23").Value = "  +0.77%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.04"
$ws.Range("E24").Value = "  +1.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.314"
$ws.Range("E25").Value = "  +1.11%  "

$ws.Range("D26").Value = "2.107.27"
$ws.Range("E26").Value = "  +1.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "159.62"
$ws.Range("E27").Value = "  +3.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.20"
$ws.Range("E28").Value = "  +1.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.110"
$ws.Range("E29").Value = "  +2.18%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.593"
$ws.Range("E30").Value = "  +2.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "121.94"
$ws.Range("E31").Value = "  +0.59%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9886"
$ws.Range("E32").Value = "  +3.94%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09498"
$ws.Range("E33").Value = "  +1.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.453"
$ws.Range("E34").Value = "  +0.39%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.617"
$ws.Range("E35").Value = "  +0.75%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.365"
$ws.Range("E36").Value = "  +1.99%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06134"
$ws.Range("E37").Value = "  +1.36%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02268"
$ws.Range("E38").Value = "  +1.79%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.231"
$ws.Range("E39").Value = "  +1.11%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.142"
$ws.Range("E40").Value = "  +0.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6014"
$ws.Range("E41").Value = "  +1.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1901"
$ws.Range("E42").Value = "  +0.65%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "10.33"
$ws.Range("E43").Value = "  +1.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.255"
$ws.Range("E44").Value = "  -2.12%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5713"
$ws.Range("E45").Value = "  +1.76%  "

$ws.Range("E46").Value = "  +1.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.949"
$ws.Range("E47").Value = "  +1.56%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.388"
$ws.Range("E48").Value = "  -0.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06914"
$ws.Range("E49").Value = "  +2.44%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "114.41"
$ws.Range("E50").Value = "  +5.50%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000303"
$ws.Range("E51").Value = "  +7.63%  "
